$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (previously "Strike#"). Update the computed K values
# for each row now that save_data regenerates K using the new formula,
# and std/mean/s_vals were recalculated accordingly.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 0
